$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 13.17295566666667
$ws.Range("N2").Value = 39.518867
$ws.Range("O2").Value = 0.133784132206724
$ws.Range("P2").Value = 0.133784132206724
$ws.Range("Q2").Value = 5.562399086851
$ws.Range("R2").Value = 50.061591781659
$ws.Range("S2").Value = 0.133784132206724
$ws.Range("T2").Value = 0.133784132206724

# Row 3
$ws.Range("O3").Value = 0.4382627974978752
$ws.Range("P3").Value = 0.4382627974978752
$ws.Range("S3").Value = 0.4382627974978752
$ws.Range("T3").Value = 0.4382627974978752

# Row 4
$ws.Range("M4").Value = 21.06166566666667
$ws.Range("N4").Value = 63.184997
$ws.Range("O4").Value = 0.2139016281041017
$ws.Range("P4").Value = 0.2139016281041017
$ws.Range("Q4").Value = 8.893477882740999
$ws.Range("R4").Value = 80.041300944669
$ws.Range("S4").Value = 0.2139016281041017
$ws.Range("T4").Value = 0.2139016281041017

# Row 5
$ws.Range("M5").Value = 21.076417
$ws.Range("N5").Value = 63.229251
$ws.Range("O5").Value = 0.214051442191299
$ws.Range("P5").Value = 0.214051442191299
$ws.Range("Q5").Value = 8.899706766003002
$ws.Range("R5").Value = 80.09736089402701
$ws.Range("S5").Value = 0.214051442191299
$ws.Range("T5").Value = 0.214051442191299
